$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("U2").Value = 1.92
$ws.Range("V2").Value = 1.77
$ws.Range("U3").Value = 1.92
$ws.Range("V3").Value = 1.77
$ws.Range("G4").Value = 1.75
$ws.Range("V4").Value = 1.63
$ws.Range("AD4").Value = 7
$ws.Range("AE4").Value = 21
$ws.Range("AO4").Value = 9.5
$ws.Range("G5").Value = 2.55
$ws.Range("I5").Value = 3.1
$ws.Range("J5").Value = 3.4
$ws.Range("V5").Value = 1.63
$ws.Range("Z5").Value = 26
$ws.Range("AI5").Value = 13
$ws.Range("AJ5").Value = 12
$ws.Range("H6").Value = 2.9
$ws.Range("I6").Value = 4.1
$ws.Range("K6").Value = 1.95
$ws.Range("AC6").Value = 6.5
$ws.Range("AH6").Value = 9
$ws.Range("AI6").Value = 17
$ws.Range("AS6").Value = 251
$ws.Range("M8").Value = 1.03
$ws.Range("O8").Value = 1.25
$ws.Range("G9").Value = 3.7
$ws.Range("H9").Value = 3.6
$ws.Range("I9").Value = 1.83
$ws.Range("J9").Value = 4.33
$ws.Range("K9").Value = 2.3
$ws.Range("M9").Value = 1.04
$ws.Range("N9").Value = 13
$ws.Range("O9").Value = 1.22
$ws.Range("AC9").Value = 13
$ws.Range("AE9").Value = 15
$ws.Range("AI9").Value = 9.5
$ws.Range("AK9").Value = 15
$ws.Range("I10").Value = 16.5
$ws.Range("J10").Value = 1.47
$ws.Range("K10").Value = 2.7
$ws.Range("L10").Value = 12
$ws.Range("AA10").Value = 9.25
$ws.Range("AL10").Value = 200
$ws.Range("AO10").Value = 4.55
$ws.Range("AQ10").Value = 10
$ws.Range("AT10").Value = 3.55
$ws.Range("AY10").Value = 70
$ws.Range("N11").Value = 15
$ws.Range("O11").Value = 1.18
$ws.Range("P11").Value = 4.5
$ws.Range("Q11").Value = 1.6
$ws.Range("R11").Value = 2.3
$ws.Range("M13").Value = 1.01
$ws.Range("N13").Value = 23
$ws.Range("AG13").Value = 81
$ws.Range("I14").Value = 1.66
$ws.Range("G15").Value = 1.72
$ws.Range("G16").Value = 1.5
$ws.Range("N16").Value = 13
$ws.Range("Q16").Value = 1.83
$ws.Range("R16").Value = 2.03
$ws.Range("Q19").Value = 1.63
$ws.Range("Q20").Value = 1.44
$ws.Range("I21").Value = 7
$ws.Range("Q21").Value = 1.3
$ws.Range("U21").Value = 1.53
$ws.Range("V21").Value = 2.38
$ws.Range("Y21").Value = 9
$ws.Range("AC21").Value = 26
$ws.Range("AK21").Value = 81
$ws.Range("AM21").Value = 41
$ws.Range("AZ21").Value = 101
$ws.Range("Q22").Value = 1.67
$ws.Range("R23").Value = 1.54
$ws.Range("M24").Value = 1.08
$ws.Range("O24").Value = 1.4
$ws.Range("R24").Value = 1.57
$ws.Range("G25").Value = 1.8
$ws.Range("M25").Value = 1.05
$ws.Range("O25").Value = 1.25
$ws.Range("M26").Value = 1.04
$ws.Range("O26").Value = 1.22
$ws.Range("Q26").Value = 1.77
$ws.Range("M27").Value = 1.08
$ws.Range("O27").Value = 1.36
$ws.Range("V27").Value = 1.69
$ws.Range("M28").Value = 1.06
$ws.Range("O28").Value = 1.29
$ws.Range("U28").Value = 1.77
$ws.Range("V28").Value = 1.87
$ws.Range("U30").Value = 1.77
$ws.Range("V30").Value = 1.87
$ws.Range("G32").Value = 2.35
$ws.Range("I32").Value = 2.75
$ws.Range("S32").Value = 1.3
$ws.Range("T32").Value = 3.4
$ws.Range("W32").Value = 11
$ws.Range("X32").Value = 13
$ws.Range("AC32").Value = 15
$ws.Range("AD32").Value = 7
$ws.Range("AE32").Value = 11
$ws.Range("AJ32").Value = 11
$ws.Range("AL32").Value = 21
$ws.Range("AT32").Value = 3.4
$ws.Range("Q33").Value = 2.25
$ws.Range("R33").Value = 1.62
$ws.Range("M35").Value = 1.05
$ws.Range("O35").Value = 1.33
$ws.Range("U35").Value = 1.87
$ws.Range("V35").Value = 1.77
$ws.Range("M36").Value = 1.03
$ws.Range("O36").Value = 1.25
$ws.Range("U36").Value = 1.69
$ws.Range("M38").Value = 1.05
$ws.Range("O38").Value = 1.41
$ws.Range("P38").Value = 2.62
$ws.Range("Q38").Value = 2.4
$ws.Range("R38").Value = 1.53
$ws.Range("V38").Value = 1.69
$ws.Range("M39").Value = 1.03
$ws.Range("O39").Value = 1.25
$ws.Range("U39").Value = 1.69
